# Estadisticos Segundo Parcial 23 Mayo
# Updates the "2o Parcial" and "Final" sheets with the grading results that
# were recorded for the second partial: counts of Aprobados/Reprobados
# (columns E/F), their percentages (G/H), the group average (I) and the
# Blancos/por_blancos columns (J/K) where applicable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "2o Parcial": rows 2-8 and 12 get real grading numbers instead of
# the placeholder "0 aprobados / 100% reprobados" values.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2o Parcial")

function Set-Row($ws, $row, $e, $f, $g, $h, $i, $j, $k) {
    $ws.Range("E" + $row).Value = $e
    $ws.Range("F" + $row).Value = $f
    $ws.Range("G" + $row).Value = $g
    $ws.Range("H" + $row).Value = $h
    $ws.Range("I" + $row).Value = $i
    $ws.Range("J" + $row).Value = $j
    $ws.Range("K" + $row).Value = $k
}

Set-Row $ws2 2  10 1  90.9 9.1  8   0  0
Set-Row $ws2 3  11 0  100  0    8.5 0  0
Set-Row $ws2 4  21 1  95.5 4.5  8.2 0  0
Set-Row $ws2 5  24 1  96   4    9.6 0  0
Set-Row $ws2 6  24 1  96   4    9.6 0  0
Set-Row $ws2 7  9  2  81.8 18.2 8   0  0
Set-Row $ws2 8  9  2  81.8 18.2 8   0  0
Set-Row $ws2 12 54 40 57.4 42.6 5.7 36 38.3

# ---------------------------------------------------------------------
# Sheet "Final": recompute the same rows using the combined 1er + 2o
# parcial figures (Blancos/por_blancos columns were already correct and
# stay untouched).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Final")

function Set-RowNoJK($ws, $row, $e, $f, $g, $h, $i) {
    $ws.Range("E" + $row).Value = $e
    $ws.Range("F" + $row).Value = $f
    $ws.Range("G" + $row).Value = $g
    $ws.Range("H" + $row).Value = $h
    $ws.Range("I" + $row).Value = $i
}

Set-RowNoJK $ws3 2  10 1 90.9 9.1 8.1
$ws3.Range("I3").Value = 8.5
Set-RowNoJK $ws3 4  21 1 95.5 4.5 8.3
$ws3.Range("I7").Value = 7.7
$ws3.Range("I8").Value = 7.7
Set-RowNoJK $ws3 12 89 5 94.7 5.3 8.2
